# Update the "results" (column G) values on the "Units" sheet of the
# Consolidated_Units reference list. These are recalculated/refreshed
# unit-conversion result values (part of the saved filter update).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# Large values are built via multiplication (with a decimal literal) to
# avoid integer-literal overflow / scientific-notation parsing issues in
# the scripting engine, while still producing an IEEE double identical to
# the target value.
$bigPercentBase = 80000000000000000000000.0   # 8E+22

$updates = [ordered]@{
    6   = 114294784
    13  = 1.91663112077723 * $bigPercentBase
    14  = 1.91663112077723 * $bigPercentBase
    23  = 1227836258968140
    25  = 1220542464
    26  = 21618063286200
    29  = 21618063286200
    45  = 58058604544
    47  = 189754293886821
    50  = 189754293886821
    51  = 58058604544
    66  = 1227836258968140
    67  = 1227836258968140
    76  = 1.91663112077723 * $bigPercentBase
    99  = 189754293886821
    100 = 883949568
    104 = 1.91663112077723 * $bigPercentBase
    105 = 1220542464
    113 = 1.91663112077723 * $bigPercentBase
    114 = 883949568
    117 = 1.91663112077723 * $bigPercentBase
    118 = 1.91663112077723 * $bigPercentBase
    120 = 1.91663112077723 * $bigPercentBase
    123 = 21618063286200
    127 = 114294784
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
